$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Joni's time-tracking hours (H) and details (I) -------------------------
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = "Setting development enviorment"

$ws.Range("H6").Value = 2
$ws.Range("I6").Value = "Creating initial sass folder and files"

$ws.Range("H7").Value = 3
$ws.Range("I7").Value = "Creating header"

$ws.Range("H10").Value = 3
$ws.Range("I10").Value = "Created 6 cards where to select functions"

$ws.Range("H11").Value = 4
$ws.Range("I11").Value = "Styled 6 cards with headings and svg icons"

$ws.Range("H13").Value = 4
$ws.Range("I13").Value = "Created popup window component"

$ws.Range("H14").Value = 3
$ws.Range("I14").Value = "Added popup to all six buttons  "

$ws.Range("H16").Value = 5
$ws.Range("I16").Value = "Styled functions inside popup"

$ws.Range("H17").Value = 5
$ws.Range("I17").Value = "Styled functions inside popup"

$ws.Range("H18").Value = 5
$ws.Range("I18").Value = "Styled functions inside popup"

$ws.Range("H19").Value = 4
$ws.Range("I19").Value = "Styled functions inside popup"

$ws.Range("G20").Value = "#########"
$ws.Range("H20").Value = 3
$ws.Range("I20").Value = "Styled functions inside popup"

$ws.Range("G21").Value = "########"
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = "Footer"

# --- Extend the total formula to include the newly used row 21 -------------
$ws.Range("H23").Formula = "=SUM(H5:H21)"

# --- Move the active selection, matching the saved workbook view -----------
$ws.Range("H10").Select()
